# "retouche de Objet definition" - update the "Objet" definition slide:
#  - add a bullet to the blank line under the "Objet" heading
#  - tweak the definition wording
#  - change "Un objet a trois..." to "On dit donc que l'objet a trois..."
#  - reposition/resize the illustrative picture

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(3)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Paragraph 2 is the empty line right after "Objet" - give it the same
# bullet (Wingdings "q") used by the "Objet" / "Class" bullets elsewhere.
$para2 = $tr.Paragraphs(2)
$bullet2 = $para2.ParagraphFormat.Bullet
$bullet2.Font.Name = "Wingdings"
$bullet2.Type = 1
$bullet2.Character = 113

# Paragraph 4: add ", un livre" before "etc."
$para4 = $tr.Paragraphs(4)
$full4 = $tr.Characters($para4.Start, $para4.Length)
$full4.Text = "Un objet est une entité qui a un état et un comportement, par exemple: une chaise, une voiture, un livre etc."

# Paragraph 5: reword the intro to the characteristics list
$para5 = $tr.Paragraphs(5)
$full5 = $tr.Characters($para5.Start, $para5.Length)
$full5.Text = "On dit donc que l'objet a trois (3) caractéristiques:"

# Reposition / resize the illustrative picture (4th shape on the slide).
# Left/Top/Width/Height are expressed in points (1 pt = 12700 EMU); the
# literals below are EMU/914400*72 nudged by <1 EMU to compensate for the
# COM layer's single-precision storage of these properties.
$pic = $s.Shapes.Item(4)
$pic.Left = 98.24708961417323
$pic.Top = 315.5850393700788
$pic.Width = 495.90882889763776
$pic.Height = 184.70173228346457
